$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value while forcing text storage (avoids numeric
# auto-conversion of values like "289.01") and preserving the cells
# original (default) style/format.
function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue "D2" '22.314.81'
$ws.Range("E2").Value = '  -5.00%  '

Set-TextValue "D3" '1.565.62'
$ws.Range("E3").Value = '  -5.13%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("E5").Value = '  +0.07%  '

Set-TextValue "D6" '289.01'
$ws.Range("E6").Value = '  -3.69%  '

Set-TextValue "D7" '0.3748'
$ws.Range("E7").Value = '  -1.02%  '

Set-TextValue "D8" '49.32'
$ws.Range("E8").Value = '  -2.58%  '

Set-TextValue "D9" '0.3418'
$ws.Range("E9").Value = '  -2.45%  '

Set-TextValue "D10" '1.165'
$ws.Range("E10").Value = '  -4.93%  '

Set-TextValue "D11" '0.07644'
$ws.Range("E11").Value = '  -5.17%  '

$ws.Range("E12").Value = '  +0.01%  '

Set-TextValue "D13" '21.40'
$ws.Range("E13").Value = '  -3.20%  '

Set-TextValue "D14" '6.009'
$ws.Range("E14").Value = '  -4.96%  '

Set-TextValue "D15" '6.939'
$ws.Range("E15").Value = '  -4.60%  '

$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue "D16" '0.00001130'
$ws.Range("E16").Value = '  -6.71%  '

$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue "D17" '1.558.20'
$ws.Range("E17").Value = '  -5.75%  '

Set-TextValue "D18" '90.01'
$ws.Range("E18").Value = '  -5.67%  '

Set-TextValue "D19" '0.06711'
$ws.Range("E19").Value = '  -3.75%  '

Set-TextValue "D20" '1.001'
$ws.Range("E20").Value = '  +0.07%  '

$ws.Range("E21").Value = '  -6.03%  '

Set-TextValue "D22" '16.58'
$ws.Range("E22").Value = '  -5.06%  '

Set-TextValue "D23" '0.5283'
$ws.Range("E23").Value = '  -8.29%  '

Set-TextValue "D24" '11.93'
$ws.Range("E24").Value = '  -4.40%  '

Set-TextValue "D25" '22.320.44'
$ws.Range("E25").Value = '  -5.00%  '

Set-TextValue "D26" '2.387'
$ws.Range("E26").Value = '  -1.54%  '

Set-TextValue "D27" '2.786'
$ws.Range("E27").Value = '  -7.61%  '

Set-TextValue "D28" '20.15'
$ws.Range("E28").Value = '  -4.48%  '

Set-TextValue "D29" '145.90'
$ws.Range("E29").Value = '  -3.87%  '

Set-TextValue "D30" '4.965'
$ws.Range("E30").Value = '  -4.25%  '

Set-TextValue "D31" '125.44'
$ws.Range("E31").Value = '  -4.88%  '

Set-TextValue "D32" '1.737.61'
$ws.Range("E32").Value = '  -5.42%  '

Set-TextValue "D33" '1.020'
$ws.Range("E33").Value = '  +3.03%  '

Set-TextValue "D34" '6.210'
$ws.Range("E34").Value = '  -10.08%  '

Set-TextValue "D35" '2.008'
$ws.Range("E35").Value = '  -6.18%  '

$ws.Range("E36").Value = '  -10.34%  '

Set-TextValue "D37" '0.08494'
$ws.Range("E37").Value = '  -3.09%  '

Set-TextValue "D38" '0.02528'
$ws.Range("E38").Value = '  -7.12%  '

Set-TextValue "D39" '0.2322'
$ws.Range("E39").Value = '  -4.38%  '

Set-TextValue "D40" '5.529'
$ws.Range("E40").Value = '  -7.00%  '

Set-TextValue "D41" '1.319'
$ws.Range("E41").Value = '  +1.79%  '

Set-TextValue "D42" '0.06392'
$ws.Range("E42").Value = '  -6.46%  '

Set-TextValue "D43" '11.72'
$ws.Range("E43").Value = '  -9.22%  '

Set-TextValue "D44" '0.6374'
$ws.Range("E44").Value = '  -7.71%  '

Set-TextValue "D45" '14.09'
$ws.Range("E45").Value = '  -9.25%  '

Set-TextValue "D46" '0.9998'
$ws.Range("E46").Value = '  +0.06%  '

Set-TextValue "D47" '0.5985'
$ws.Range("E47").Value = '  -6.45%  '

Set-TextValue "D48" '3.751'
$ws.Range("E48").Value = '  -4.29%  '

$ws.Range("E49").Value = '  -7.05%  '

Set-TextValue "D50" '1.268'
$ws.Range("E50").Value = '  +2.23%  '

Set-TextValue "D51" '124.14'
$ws.Range("E51").Value = '  -2.41%  '
